$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap data between rows 8 and 9 (id/Div/Date columns A, C, D stay put)
$ws.Cells.Item(8, 2).Value = 6228028
$ws.Cells.Item(9, 2).Value = 6228588
$ws.Cells.Item(8, 5).Value = "Lion City Sailors FC"
$ws.Cells.Item(9, 5).Value = "Tampines Rovers FC"
$ws.Cells.Item(8, 6).Value = "Young Lions"
$ws.Cells.Item(9, 6).Value = "Hougang United FC"
$ws.Cells.Item(8, 7).Value = 4
$ws.Cells.Item(9, 7).Value = 2
$ws.Cells.Item(8, 8).Value = 1
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(8, 9).Value = "H"
$ws.Cells.Item(9, 9).Value = "H"
$ws.Cells.Item(8, 10).Value = 1.083
$ws.Cells.Item(9, 10).Value = 1.333
$ws.Cells.Item(8, 11).Value = 9
$ws.Cells.Item(9, 11).Value = 5
$ws.Cells.Item(8, 12).Value = 15
$ws.Cells.Item(9, 12).Value = 6.25
$ws.Cells.Item(8, 13).Value = 1.142
$ws.Cells.Item(9, 13).Value = 1.615
$ws.Cells.Item(8, 14).Value = 7
$ws.Cells.Item(9, 14).Value = 4.5
$ws.Cells.Item(8, 15).Value = 12
$ws.Cells.Item(9, 15).Value = 3.8
$ws.Cells.Item(8, 16).Value = -2.25
$ws.Cells.Item(9, 16).Value = -1
$ws.Cells.Item(8, 17).Value = 1.85
$ws.Cells.Item(9, 17).Value = 2.05
$ws.Cells.Item(8, 18).Value = 2
$ws.Cells.Item(9, 18).Value = 1.8
$ws.Cells.Item(8, 19).Value = 4
$ws.Cells.Item(9, 19).Value = 3.5
$ws.Cells.Item(8, 20).Value = 1.875
$ws.Cells.Item(9, 20).Value = 1.85
$ws.Cells.Item(8, 21).Value = 1.975
$ws.Cells.Item(9, 21).Value = 2
$ws.Cells.Item(8, 22).Value = 0.1419999999999999
$ws.Cells.Item(9, 22).Value = 0.615
$ws.Cells.Item(8, 23).Value = -1
$ws.Cells.Item(9, 23).Value = -1
$ws.Cells.Item(8, 24).Value = -1
$ws.Cells.Item(9, 24).Value = -1
$ws.Cells.Item(8, 25).Value = 0.8500000000000001
$ws.Cells.Item(9, 25).Value = 1.05
$ws.Cells.Item(8, 26).Value = -1
$ws.Cells.Item(9, 26).Value = -1
$ws.Cells.Item(8, 27).Value = 0.875
$ws.Cells.Item(9, 27).Value = -1
$ws.Cells.Item(8, 28).Value = -1
$ws.Cells.Item(9, 28).Value = 1

# Swap data between rows 13 and 14 (id/Div/Date columns A, C, D stay put)
$ws.Cells.Item(13, 2).Value = 6228592
$ws.Cells.Item(14, 2).Value = 6228593
$ws.Cells.Item(13, 5).Value = "Albirex Niigata Singapore"
$ws.Cells.Item(14, 5).Value = "DPMM FC"
$ws.Cells.Item(13, 6).Value = "Tanjong Pagar United"
$ws.Cells.Item(14, 6).Value = "Young Lions"
$ws.Cells.Item(13, 7).Value = 4
$ws.Cells.Item(14, 7).Value = 6
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(13, 9).Value = "H"
$ws.Cells.Item(14, 9).Value = "H"
$ws.Cells.Item(13, 10).Value = 1.166
$ws.Cells.Item(14, 10).Value = 1.444
$ws.Cells.Item(13, 11).Value = 7
$ws.Cells.Item(14, 11).Value = 4.5
$ws.Cells.Item(13, 12).Value = 9
$ws.Cells.Item(14, 12).Value = 5
$ws.Cells.Item(13, 13).Value = 1.125
$ws.Cells.Item(14, 13).Value = 1.444
$ws.Cells.Item(13, 14).Value = 8
$ws.Cells.Item(14, 14).Value = 4.5
$ws.Cells.Item(13, 15).Value = 15
$ws.Cells.Item(14, 15).Value = 5
$ws.Cells.Item(13, 16).Value = -2.5
$ws.Cells.Item(14, 16).Value = -1.25
$ws.Cells.Item(13, 17).Value = 1.95
$ws.Cells.Item(14, 17).Value = 2
$ws.Cells.Item(13, 18).Value = 1.9
$ws.Cells.Item(14, 18).Value = 1.85
$ws.Cells.Item(13, 19).Value = 4
$ws.Cells.Item(14, 19).Value = 3.25
$ws.Cells.Item(13, 20).Value = 1.925
$ws.Cells.Item(14, 20).Value = 1.875
$ws.Cells.Item(13, 21).Value = 1.925
$ws.Cells.Item(14, 21).Value = 1.975
$ws.Cells.Item(13, 22).Value = 0.125
$ws.Cells.Item(14, 22).Value = 0.444
$ws.Cells.Item(13, 23).Value = -1
$ws.Cells.Item(14, 23).Value = -1
$ws.Cells.Item(13, 24).Value = -1
$ws.Cells.Item(14, 24).Value = -1
$ws.Cells.Item(13, 25).Value = 0.95
$ws.Cells.Item(14, 25).Value = 1
$ws.Cells.Item(13, 26).Value = -1
$ws.Cells.Item(14, 26).Value = -1
$ws.Cells.Item(13, 27).Value = 0
$ws.Cells.Item(14, 27).Value = 0.875
$ws.Cells.Item(13, 28).Value = 0
$ws.Cells.Item(14, 28).Value = -1

# Swap data between rows 20 and 21 (id/Div/Date columns A, C, D stay put)
$ws.Cells.Item(20, 2).Value = 6228597
$ws.Cells.Item(21, 2).Value = 6228598
$ws.Cells.Item(20, 5).Value = "Hougang United FC"
$ws.Cells.Item(21, 5).Value = "Geylang International"
$ws.Cells.Item(20, 6).Value = "Balestier Khalsa FC"
$ws.Cells.Item(21, 6).Value = "Tampines Rovers FC"
$ws.Cells.Item(20, 7).Value = 1
$ws.Cells.Item(21, 7).Value = 1
$ws.Cells.Item(20, 8).Value = 3
$ws.Cells.Item(21, 8).Value = 1
$ws.Cells.Item(20, 9).Value = "A"
$ws.Cells.Item(21, 9).Value = "D"
$ws.Cells.Item(20, 10).Value = 2.5
$ws.Cells.Item(21, 10).Value = 3.6
$ws.Cells.Item(20, 11).Value = 3.6
$ws.Cells.Item(21, 11).Value = 4.2
$ws.Cells.Item(20, 12).Value = 2.25
$ws.Cells.Item(21, 12).Value = 1.666
$ws.Cells.Item(20, 13).Value = 2.6
$ws.Cells.Item(21, 13).Value = 4
$ws.Cells.Item(20, 14).Value = 3.75
$ws.Cells.Item(21, 14).Value = 4.5
$ws.Cells.Item(20, 15).Value = 2.2
$ws.Cells.Item(21, 15).Value = 1.55
$ws.Cells.Item(20, 16).Value = 0.25
$ws.Cells.Item(21, 16).Value = 1
$ws.Cells.Item(20, 17).Value = 1.825
$ws.Cells.Item(21, 17).Value = 1.85
$ws.Cells.Item(20, 18).Value = 2.025
$ws.Cells.Item(21, 18).Value = 2
$ws.Cells.Item(20, 19).Value = 4
$ws.Cells.Item(21, 19).Value = 3.25
$ws.Cells.Item(20, 20).Value = 1.95
$ws.Cells.Item(21, 20).Value = 2.025
$ws.Cells.Item(20, 21).Value = 1.9
$ws.Cells.Item(21, 21).Value = 1.825
$ws.Cells.Item(20, 22).Value = -1
$ws.Cells.Item(21, 22).Value = -1
$ws.Cells.Item(20, 23).Value = -1
$ws.Cells.Item(21, 23).Value = 3.5
$ws.Cells.Item(20, 24).Value = 1.2
$ws.Cells.Item(21, 24).Value = -1
$ws.Cells.Item(20, 25).Value = -1
$ws.Cells.Item(21, 25).Value = 0.8500000000000001
$ws.Cells.Item(20, 26).Value = 1.025
$ws.Cells.Item(21, 26).Value = -1
$ws.Cells.Item(20, 27).Value = 0
$ws.Cells.Item(21, 27).Value = -1
$ws.Cells.Item(20, 28).Value = 0
$ws.Cells.Item(21, 28).Value = 0.825

# Swap data between rows 26 and 27 (id/Div/Date columns A, C, D stay put)
$ws.Cells.Item(26, 2).Value = 6228032
$ws.Cells.Item(27, 2).Value = 6228602
$ws.Cells.Item(26, 5).Value = "Lion City Sailors FC"
$ws.Cells.Item(27, 5).Value = "Tampines Rovers FC"
$ws.Cells.Item(26, 6).Value = "DPMM FC"
$ws.Cells.Item(27, 6).Value = "Tanjong Pagar United"
$ws.Cells.Item(26, 7).Value = 1
$ws.Cells.Item(27, 7).Value = 2
$ws.Cells.Item(26, 8).Value = 3
$ws.Cells.Item(27, 8).Value = 1
$ws.Cells.Item(26, 9).Value = "A"
$ws.Cells.Item(27, 9).Value = "H"
$ws.Cells.Item(26, 10).Value = 1.2
$ws.Cells.Item(27, 10).Value = 1.2
$ws.Cells.Item(26, 11).Value = 6
$ws.Cells.Item(27, 11).Value = 6
$ws.Cells.Item(26, 12).Value = 9
$ws.Cells.Item(27, 12).Value = 9
$ws.Cells.Item(26, 13).Value = 1.142
$ws.Cells.Item(27, 13).Value = 1.125
$ws.Cells.Item(26, 14).Value = 7.5
$ws.Cells.Item(27, 14).Value = 7.5
$ws.Cells.Item(26, 15).Value = 13
$ws.Cells.Item(27, 15).Value = 17
$ws.Cells.Item(26, 16).Value = -2.25
$ws.Cells.Item(27, 16).Value = -2.5
$ws.Cells.Item(26, 17).Value = 1.825
$ws.Cells.Item(27, 17).Value = 1.9
$ws.Cells.Item(26, 18).Value = 2.025
$ws.Cells.Item(27, 18).Value = 1.95
$ws.Cells.Item(26, 19).Value = 4
$ws.Cells.Item(27, 19).Value = 4.25
$ws.Cells.Item(26, 20).Value = 1.825
$ws.Cells.Item(27, 20).Value = 1.975
$ws.Cells.Item(26, 21).Value = 2.025
$ws.Cells.Item(27, 21).Value = 1.875
$ws.Cells.Item(26, 22).Value = -1
$ws.Cells.Item(27, 22).Value = 0.125
$ws.Cells.Item(26, 23).Value = -1
$ws.Cells.Item(27, 23).Value = -1
$ws.Cells.Item(26, 24).Value = 12
$ws.Cells.Item(27, 24).Value = -1
$ws.Cells.Item(26, 25).Value = -1
$ws.Cells.Item(27, 25).Value = -1
$ws.Cells.Item(26, 26).Value = 1.025
$ws.Cells.Item(27, 26).Value = 0.95
$ws.Cells.Item(26, 27).Value = 0
$ws.Cells.Item(27, 27).Value = -1
$ws.Cells.Item(26, 28).Value = 0
$ws.Cells.Item(27, 28).Value = 0.875

# Swap data between rows 32 and 33 (id/Div/Date columns A, C, D stay put)
$ws.Cells.Item(32, 2).Value = 6915906
$ws.Cells.Item(33, 2).Value = 6228033
$ws.Cells.Item(32, 5).Value = "DPMM FC"
$ws.Cells.Item(33, 5).Value = "Albirex Niigata Singapore"
$ws.Cells.Item(32, 6).Value = "Balestier Khalsa FC"
$ws.Cells.Item(33, 6).Value = "Lion City Sailors FC"
$ws.Cells.Item(32, 7).Value = 2
$ws.Cells.Item(33, 7).Value = 3
$ws.Cells.Item(32, 8).Value = 3
$ws.Cells.Item(33, 8).Value = 1
$ws.Cells.Item(32, 9).Value = "A"
$ws.Cells.Item(33, 9).Value = "H"
$ws.Cells.Item(32, 10).Value = 2.1
$ws.Cells.Item(33, 10).Value = 1.6
$ws.Cells.Item(32, 11).Value = 4
$ws.Cells.Item(33, 11).Value = 4.5
$ws.Cells.Item(32, 12).Value = 2.6
$ws.Cells.Item(33, 12).Value = 3.75
$ws.Cells.Item(32, 13).Value = 2.05
$ws.Cells.Item(33, 13).Value = 1.8
$ws.Cells.Item(32, 14).Value = 4.5
$ws.Cells.Item(33, 14).Value = 4.2
$ws.Cells.Item(32, 15).Value = 2.5
$ws.Cells.Item(33, 15).Value = 3.1
$ws.Cells.Item(32, 16).Value = -0.25
$ws.Cells.Item(33, 16).Value = -0.5
$ws.Cells.Item(32, 17).Value = 1.925
$ws.Cells.Item(33, 17).Value = 1.85
$ws.Cells.Item(32, 18).Value = 1.925
$ws.Cells.Item(33, 18).Value = 2
$ws.Cells.Item(32, 19).Value = 4
$ws.Cells.Item(33, 19).Value = 3.75
$ws.Cells.Item(32, 20).Value = 1.9
$ws.Cells.Item(33, 20).Value = 1.925
$ws.Cells.Item(32, 21).Value = 1.95
$ws.Cells.Item(33, 21).Value = 1.925
$ws.Cells.Item(32, 22).Value = -1
$ws.Cells.Item(33, 22).Value = 0.8
$ws.Cells.Item(32, 23).Value = -1
$ws.Cells.Item(33, 23).Value = -1
$ws.Cells.Item(32, 24).Value = 1.5
$ws.Cells.Item(33, 24).Value = -1
$ws.Cells.Item(32, 25).Value = -1
$ws.Cells.Item(33, 25).Value = 0.8500000000000001
$ws.Cells.Item(32, 26).Value = 0.925
$ws.Cells.Item(33, 26).Value = -1
$ws.Cells.Item(32, 27).Value = 0.8999999999999999
$ws.Cells.Item(33, 27).Value = 0.4625
$ws.Cells.Item(32, 28).Value = -1
$ws.Cells.Item(33, 28).Value = -0.5

# Swap data between rows 36 and 37 (id/Div/Date columns A, C, D stay put)
$ws.Cells.Item(36, 2).Value = 6228609
$ws.Cells.Item(37, 2).Value = 6228610
$ws.Cells.Item(36, 5).Value = "Albirex Niigata Singapore"
$ws.Cells.Item(37, 5).Value = "Tampines Rovers FC"
$ws.Cells.Item(36, 6).Value = "Young Lions"
$ws.Cells.Item(37, 6).Value = "Geylang International"
$ws.Cells.Item(36, 7).Value = 5
$ws.Cells.Item(37, 7).Value = 2
$ws.Cells.Item(36, 8).Value = 0
$ws.Cells.Item(37, 8).Value = 3
$ws.Cells.Item(36, 9).Value = "H"
$ws.Cells.Item(37, 9).Value = "A"
$ws.Cells.Item(36, 10).Value = 1.062
$ws.Cells.Item(37, 10).Value = 1.333
$ws.Cells.Item(36, 11).Value = 11
$ws.Cells.Item(37, 11).Value = 4.5
$ws.Cells.Item(36, 12).Value = 17
$ws.Cells.Item(37, 12).Value = 7
$ws.Cells.Item(36, 13).Value = 1.055
$ws.Cells.Item(37, 13).Value = 1.25
$ws.Cells.Item(36, 14).Value = 11
$ws.Cells.Item(37, 14).Value = 5.5
$ws.Cells.Item(36, 15).Value = 17
$ws.Cells.Item(37, 15).Value = 8
$ws.Cells.Item(36, 16).Value = -3.25
$ws.Cells.Item(37, 16).Value = -1.75
$ws.Cells.Item(36, 17).Value = 1.925
$ws.Cells.Item(37, 17).Value = 1.9
$ws.Cells.Item(36, 18).Value = 1.925
$ws.Cells.Item(37, 18).Value = 1.95
$ws.Cells.Item(36, 19).Value = 4.5
$ws.Cells.Item(37, 19).Value = 3.75
$ws.Cells.Item(36, 20).Value = 1.85
$ws.Cells.Item(37, 20).Value = 1.925
$ws.Cells.Item(36, 21).Value = 2
$ws.Cells.Item(37, 21).Value = 1.925
$ws.Cells.Item(36, 22).Value = 0.05499999999999994
$ws.Cells.Item(37, 22).Value = -1
$ws.Cells.Item(36, 23).Value = -1
$ws.Cells.Item(37, 23).Value = -1
$ws.Cells.Item(36, 24).Value = -1
$ws.Cells.Item(37, 24).Value = 7
$ws.Cells.Item(36, 25).Value = 0.925
$ws.Cells.Item(37, 25).Value = -1
$ws.Cells.Item(36, 26).Value = -1
$ws.Cells.Item(37, 26).Value = 0.95
$ws.Cells.Item(36, 27).Value = 0.8500000000000001
$ws.Cells.Item(37, 27).Value = 0.925
$ws.Cells.Item(36, 28).Value = -1
$ws.Cells.Item(37, 28).Value = -1

# Swap data between rows 43 and 44 (id/Div/Date columns A, C, D stay put)
$ws.Cells.Item(43, 2).Value = 6228036
$ws.Cells.Item(44, 2).Value = 6228615
$ws.Cells.Item(43, 5).Value = "Tampines Rovers FC"
$ws.Cells.Item(44, 5).Value = "Hougang United FC"
$ws.Cells.Item(43, 6).Value = "Lion City Sailors FC"
$ws.Cells.Item(44, 6).Value = "DPMM FC"
$ws.Cells.Item(43, 7).Value = 2
$ws.Cells.Item(44, 7).Value = 1
$ws.Cells.Item(43, 8).Value = 5
$ws.Cells.Item(44, 8).Value = 0
$ws.Cells.Item(43, 9).Value = "A"
$ws.Cells.Item(44, 9).Value = "H"
$ws.Cells.Item(43, 10).Value = 2.3
$ws.Cells.Item(44, 10).Value = 1.909
$ws.Cells.Item(43, 11).Value = 4.2
$ws.Cells.Item(44, 11).Value = 4.5
$ws.Cells.Item(43, 12).Value = 2.25
$ws.Cells.Item(44, 12).Value = 2.75
$ws.Cells.Item(43, 13).Value = 2.5
$ws.Cells.Item(44, 13).Value = 1.909
$ws.Cells.Item(43, 14).Value = 4.2
$ws.Cells.Item(44, 14).Value = 4.5
$ws.Cells.Item(43, 15).Value = 2.1
$ws.Cells.Item(44, 15).Value = 2.7
$ws.Cells.Item(43, 16).Value = 0.25
$ws.Cells.Item(44, 16).Value = -0.5
$ws.Cells.Item(43, 17).Value = 1.875
$ws.Cells.Item(44, 17).Value = 2
$ws.Cells.Item(43, 18).Value = 1.975
$ws.Cells.Item(44, 18).Value = 1.85
$ws.Cells.Item(43, 19).Value = 4
$ws.Cells.Item(44, 19).Value = 4
$ws.Cells.Item(43, 20).Value = 1.95
$ws.Cells.Item(44, 20).Value = 2
$ws.Cells.Item(43, 21).Value = 1.9
$ws.Cells.Item(44, 21).Value = 1.85
$ws.Cells.Item(43, 22).Value = -1
$ws.Cells.Item(44, 22).Value = 0.909
$ws.Cells.Item(43, 23).Value = -1
$ws.Cells.Item(44, 23).Value = -1
$ws.Cells.Item(43, 24).Value = 1.1
$ws.Cells.Item(44, 24).Value = -1
$ws.Cells.Item(43, 25).Value = -1
$ws.Cells.Item(44, 25).Value = 1
$ws.Cells.Item(43, 26).Value = 0.9750000000000001
$ws.Cells.Item(44, 26).Value = -1
$ws.Cells.Item(43, 27).Value = 0.95
$ws.Cells.Item(44, 27).Value = -1
$ws.Cells.Item(43, 28).Value = -1
$ws.Cells.Item(44, 28).Value = 0.8500000000000001

# Swap data between rows 47 and 48 (id/Div/Date columns A, C, D stay put)
$ws.Cells.Item(47, 2).Value = 6228618
$ws.Cells.Item(48, 2).Value = 6228619
$ws.Cells.Item(47, 5).Value = "Young Lions"
$ws.Cells.Item(48, 5).Value = "Geylang International"
$ws.Cells.Item(47, 6).Value = "DPMM FC"
$ws.Cells.Item(48, 6).Value = "Balestier Khalsa FC"
$ws.Cells.Item(47, 7).Value = 0
$ws.Cells.Item(48, 7).Value = 2
$ws.Cells.Item(47, 8).Value = 2
$ws.Cells.Item(48, 8).Value = 6
$ws.Cells.Item(47, 9).Value = "A"
$ws.Cells.Item(48, 9).Value = "A"
$ws.Cells.Item(47, 10).Value = 5
$ws.Cells.Item(48, 10).Value = 2.1
$ws.Cells.Item(47, 11).Value = 5
$ws.Cells.Item(48, 11).Value = 4
$ws.Cells.Item(47, 12).Value = 1.4
$ws.Cells.Item(48, 12).Value = 2.55
$ws.Cells.Item(47, 13).Value = 5
$ws.Cells.Item(48, 13).Value = 2.15
$ws.Cells.Item(47, 14).Value = 5.25
$ws.Cells.Item(48, 14).Value = 3.6
$ws.Cells.Item(47, 15).Value = 1.4
$ws.Cells.Item(48, 15).Value = 2.8
$ws.Cells.Item(47, 16).Value = 1.25
$ws.Cells.Item(48, 16).Value = -0.25
$ws.Cells.Item(47, 17).Value = 2
$ws.Cells.Item(48, 17).Value = 1.95
$ws.Cells.Item(47, 18).Value = 1.85
$ws.Cells.Item(48, 18).Value = 1.9
$ws.Cells.Item(47, 19).Value = 4
$ws.Cells.Item(48, 19).Value = 4.5
$ws.Cells.Item(47, 20).Value = 1.975
$ws.Cells.Item(48, 20).Value = 2
$ws.Cells.Item(47, 21).Value = 1.875
$ws.Cells.Item(48, 21).Value = 1.85
$ws.Cells.Item(47, 22).Value = -1
$ws.Cells.Item(48, 22).Value = -1
$ws.Cells.Item(47, 23).Value = -1
$ws.Cells.Item(48, 23).Value = -1
$ws.Cells.Item(47, 24).Value = 0.3999999999999999
$ws.Cells.Item(48, 24).Value = 1.8
$ws.Cells.Item(47, 25).Value = -1
$ws.Cells.Item(48, 25).Value = -1
$ws.Cells.Item(47, 26).Value = 0.8500000000000001
$ws.Cells.Item(48, 26).Value = 0.8999999999999999
$ws.Cells.Item(47, 27).Value = -1
$ws.Cells.Item(48, 27).Value = 1
$ws.Cells.Item(47, 28).Value = 0.875
$ws.Cells.Item(48, 28).Value = -1

# Swap data between rows 51 and 52 (id/Div/Date columns A, C, D stay put)
$ws.Cells.Item(51, 2).Value = 6228622
$ws.Cells.Item(52, 2).Value = 6228621
$ws.Cells.Item(51, 5).Value = "Hougang United FC"
$ws.Cells.Item(52, 5).Value = "Geylang International"
$ws.Cells.Item(51, 6).Value = "Tanjong Pagar United"
$ws.Cells.Item(52, 6).Value = "Young Lions"
$ws.Cells.Item(51, 7).Value = 3
$ws.Cells.Item(52, 7).Value = 3
$ws.Cells.Item(51, 8).Value = 3
$ws.Cells.Item(52, 8).Value = 0
$ws.Cells.Item(51, 9).Value = "D"
$ws.Cells.Item(52, 9).Value = "H"
$ws.Cells.Item(51, 10).Value = 1.727
$ws.Cells.Item(52, 10).Value = 1.363
$ws.Cells.Item(51, 11).Value = 4.2
$ws.Cells.Item(52, 11).Value = 5
$ws.Cells.Item(51, 12).Value = 3.4
$ws.Cells.Item(52, 12).Value = 5.75
$ws.Cells.Item(51, 13).Value = 1.727
$ws.Cells.Item(52, 13).Value = 1.285
$ws.Cells.Item(51, 14).Value = 4.333
$ws.Cells.Item(52, 14).Value = 6
$ws.Cells.Item(51, 15).Value = 3.5
$ws.Cells.Item(52, 15).Value = 6.5
$ws.Cells.Item(51, 16).Value = -0.75
$ws.Cells.Item(52, 16).Value = -1.75
$ws.Cells.Item(51, 17).Value = 1.975
$ws.Cells.Item(52, 17).Value = 1.975
$ws.Cells.Item(51, 18).Value = 1.875
$ws.Cells.Item(52, 18).Value = 1.875
$ws.Cells.Item(51, 19).Value = 4.25
$ws.Cells.Item(52, 19).Value = 4
$ws.Cells.Item(51, 20).Value = 2
$ws.Cells.Item(52, 20).Value = 1.825
$ws.Cells.Item(51, 21).Value = 1.85
$ws.Cells.Item(52, 21).Value = 2.025
$ws.Cells.Item(51, 22).Value = -1
$ws.Cells.Item(52, 22).Value = 0.2849999999999999
$ws.Cells.Item(51, 23).Value = 3.333
$ws.Cells.Item(52, 23).Value = -1
$ws.Cells.Item(51, 24).Value = -1
$ws.Cells.Item(52, 24).Value = -1
$ws.Cells.Item(51, 25).Value = -1
$ws.Cells.Item(52, 25).Value = 0.9750000000000001
$ws.Cells.Item(51, 26).Value = 0.875
$ws.Cells.Item(52, 26).Value = -1
$ws.Cells.Item(51, 27).Value = 1
$ws.Cells.Item(52, 27).Value = -1
$ws.Cells.Item(51, 28).Value = -1
$ws.Cells.Item(52, 28).Value = 1.025
